# Automatic update of files.
# Rows 2-9 of the "Artfynd" sheet get their species-occurrence data
# reshuffled: each row's Id/Taxon/RedList/TaxonId/Name/SciName/Author/
# East/North/Comment fields are replaced with another record's values
# (the AC "Publik kommentar" comment cell is removed from row 5 and a
# new one is added to row 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111639169
$ws.Range("Q2").Value = 548224.5774945696
$ws.Range("R2").Value = 6926512.579557057
$ws.Range("AC2").Value = 'riklig förekomst, mer än 50 plantor'

# Row 3
$ws.Range("A3").Value = 111639172
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = 'VU'
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = 'Knärot'
$ws.Range("G3").Value = 'Goodyera repens'
$ws.Range("H3").Value = '(L.) R. Br.'
$ws.Range("Q3").Value = 548221.3480213688
$ws.Range("R3").Value = 6926511.607424877

# Row 4
$ws.Range("A4").Value = 111639174
$ws.Range("Q4").Value = 547803.9854679118
$ws.Range("R4").Value = 6926147.447742103
$ws.Range("AC4").Value = 'ca 6 plantor'

# Row 5
$ws.Range("A5").Value = 111639175
$ws.Range("B5").Value = 89686
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 658
$ws.Range("F5").Value = 'Rosenticka'
$ws.Range("G5").Value = 'Rhodofomes roseus'
$ws.Range("H5").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q5").Value = 547828.4099300706
$ws.Range("R5").Value = 6926124.660841302
$ws.Range("AC5").ClearContents()

# Row 6
$ws.Range("A6").Value = 111639168
$ws.Range("B6").Value = 89686
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 658
$ws.Range("F6").Value = 'Rosenticka'
$ws.Range("G6").Value = 'Rhodofomes roseus'
$ws.Range("H6").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q6").Value = 548104.1391889038
$ws.Range("R6").Value = 6926477.987023209

# Row 7
$ws.Range("A7").Value = 111639170
$ws.Range("Q7").Value = 548231.4260436196
$ws.Range("R7").Value = 6926519.619127685
$ws.Range("AC7").Value = 'ca 15 plantor'

# Row 8
$ws.Range("A8").Value = 111639173
$ws.Range("Q8").Value = 547838.0352795018
$ws.Range("R8").Value = 6926228.915831603
$ws.Range("AC8").Value = 'ca 15 plantor'

# Row 9
$ws.Range("A9").Value = 111639167
$ws.Range("B9").Value = 96348
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = 'Knärot'
$ws.Range("G9").Value = 'Goodyera repens'
$ws.Range("H9").Value = '(L.) R. Br.'
$ws.Range("Q9").Value = 547814.5103353403
$ws.Range("R9").Value = 6926124.461383951
$ws.Range("AC9").Value = '1 planta'
